# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.978.66'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '''1.897.21'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''0.8284'
$ws.Range("E5").Value = '  +4.53%  '
$ws.Range("D6").Value = '''241.74'
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("D7").Value = '''1.001'
$ws.Range("D8").Value = '''0.3269'
$ws.Range("E8").Value = '  +2.69%  '
$ws.Range("D9").Value = '''26.44'
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = '''0.07022'
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("D11").Value = '''0.08088'
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D12").Value = '''0.7596'
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").Value = '''1.898.76'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").Value = '''5.238'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("D16").Value = '''29.983.74'
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").Value = '''14.07'
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("D18").Value = '''5.839'
$ws.Range("E18").Value = '  -2.71%  '
$ws.Range("D19").Value = '''243.45'
$ws.Range("E19").Value = '  -2.30%  '
$ws.Range("D20").Value = '''0.000007741'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").Value = '''1.001'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").Value = '''2.149.76'
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '''6.941'
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("D25").Value = '''0.1719'
$ws.Range("E25").Value = '  +22.87%  '
$ws.Range("D26").Value = '''9.242'
$ws.Range("E26").Value = '  -1.01%  '
$ws.Range("D27").Value = '''165.43'
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").Value = '''18.87'
$ws.Range("E28").Value = '  -0.75%  '
$ws.Range("D29").Value = '''2.088'
$ws.Range("E29").Value = '  +1.54%  '
$ws.Range("D30").Value = '''1.360'
$ws.Range("E30").Value = '  -2.29%  '
$ws.Range("D31").Value = '''1.511'
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").Value = '''0.05911'
$ws.Range("E32").Value = '  +9.22%  '
$ws.Range("D33").Value = '''4.269'
$ws.Range("E33").Value = '  -2.33%  '
$ws.Range("D34").Value = '''4.059'
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("D35").Value = '''1.263'
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("D36").Value = '''0.7295'
$ws.Range("D37").Value = '''2.721'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").Value = '''0.01912'
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("D39").Value = '''2.776'
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("D40").Value = '''0.4430'
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").Value = '''72.30'
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("D42").Value = '''5.846'
$ws.Range("E42").Value = '  -5.55%  '
$ws.Range("D43").Value = '''0.8508'
$ws.Range("E43").Value = '  +1.91%  '
$ws.Range("D44").Value = '''1.001'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '''1.893'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("D46").Value = '''101.94'
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("D47").Value = '''7.533'
$ws.Range("E47").Value = '  -1.13%  '
$ws.Range("D48").Value = '''9.755'
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("D49").Value = '''988.15'
$ws.Range("E49").Value = '  +2.47%  '
$ws.Range("D50").Value = '''2.047.76'
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").Value = '''1.514'
$ws.Range("E51").Value = '  +0.06%  '
